$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2016")

# Update raw input values (D = price per share, R = Dec dividend column)
$ws.Range("D4").Value = 20.352
$ws.Range("R4").Value = 9.69

$ws.Range("D7").Value = 10.07
$ws.Range("R7").Value = 5.96

$ws.Range("D8").Value = 8.124
$ws.Range("R8").Value = 7.54

# Match the saved cursor position from the edited workbook
$ws.Activate() | Out-Null
$ws.Range("I15").Select() | Out-Null
